# Updated data model: populate the tr0001 (test run) summary sheet with
# actual uuid / build / test-suite / test-case values, and refresh the
# active sheet / selected cell bookkeeping to match the new editing state.

$wb = $excel.ActiveWorkbook

# --- tr0001: fill in the new B-column values and the new D-column block ---
$ws1 = $wb.Worksheets.Item("tr0001")

$ws1.Range("A1").Value = "tr_uuid:"
$ws1.Range("B1").Value = "tr100101"
$ws1.Range("C1").Value = "start time:"

$ws1.Range("A2").Value = "build:"
$ws1.Range("B2").Value = "V2.4"
$ws1.Range("C2").Value = "end time:"

$ws1.Range("A3").Value = "test suites:"
$ws1.Range("B3").Value = "uuid-ts-01"
$ws1.Range("C3").Value = "test cases:"
$ws1.Range("D3").Value = "uuid-tc-01"

$ws1.Range("D4").Value = "Uuid-tc-02"
$ws1.Range("D5").Value = "dfj83jf-hf84-kg83-jhsh3"
$ws1.Range("D6").Value = "jf37gu-if83jf-8utgh4-gjr84"

# --- refresh each sheet's remembered selection ---
$ws2 = $wb.Worksheets.Item("ts0001")
$ws2.Range("B1").Select()

$ws3 = $wb.Worksheets.Item("tc0001")
$ws3.Range("B1").Select()

$ws4 = $wb.Worksheets.Item("tc0002")
$ws4.Range("B1").Select()

$ws5 = $wb.Worksheets.Item("asdf")
$ws5.Range("B1").Select()

$ws6 = $wb.Worksheets.Item("827asaf")
$ws6.Range("B1").Select()

# tr0001 becomes the active sheet again, with C13 selected
$ws1.Activate()
$ws1.Range("C13").Select()
